$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted right before the current
# row 128 ("Jengibre" / Vega Central Mapocho de Santiago), pushing the
# former rows 128-136 down to 129-137.
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new observation's data.
$ws.Cells.Item(128, 1).Value = 9
$ws.Cells.Item(128, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(128, 3).Value = "Metropolitana"
$ws.Cells.Item(128, 4).Value = 45124
$ws.Cells.Item(128, 5).Value = 13
$ws.Cells.Item(128, 6).Value = 100114007
$ws.Cells.Item(128, 7).Value = "Jengibre"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 430
$ws.Cells.Item(128, 11).Value = 15000
$ws.Cells.Item(128, 12).Value = 16000
$ws.Cells.Item(128, 13).Value = 15500
$ws.Cells.Item(128, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(128, 15).Value = "Perú"
$ws.Cells.Item(128, 16).Value = 1192
$ws.Cells.Item(128, 17).Value = 13
$ws.Cells.Item(128, 18).Value = "Hortaliza"
